$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (ALC)
$ws.Range("H18").Value = 1187.0555
$ws.Range("I18").Value = 965.25
$ws.Range("J18").Value = 1630.6666
$ws.Range("K18").Value = 965.25
$ws.Range("L18").Value = 1630.6666
$ws.Range("M18").Value = -681.25
$ws.Range("N18").Value = -2198.6666

# Row 40 (ALC)
$ws.Range("H40").Value = 1814.1
$ws.Range("I40").Value = 1571.2222
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 1571.2222
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1396.2222
$ws.Range("N40").Value = -4350

# Row 100 (ALC)
$ws.Range("H100").Value = 49385610
$ws.Range("I100").Value = 88890770
$ws.Range("J100").Value = 4150
$ws.Range("K100").Value = 88890770
$ws.Range("L100").Value = 4150
$ws.Range("M100").Value = -88890229
$ws.Range("N100").Value = -5232

# Row 137 (ALC)
$ws.Range("H137").Value = 1693.4375
$ws.Range("I137").Value = 1174.5834
$ws.Range("J137").Value = 3250
$ws.Range("K137").Value = 3523.7502
$ws.Range("L137").Value = 9750
$ws.Range("M137").Value = -973.7501999999999
$ws.Range("N137").Value = -14850

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 6376.5557
$ws.Range("I32").Value = 3205.541
$ws.Range("J32").Value = 23961.273
$ws.Range("K32").Value = 3205.541
$ws.Range("L32").Value = 23961.273
$ws.Range("M32").Value = -2918.541
$ws.Range("N32").Value = -24535.273

# Row 102 (ARM)
$ws.Range("H102").Value = 1651
$ws.Range("I102").Value = 1610.5264
$ws.Range("J102").Value = 1907.3334
$ws.Range("K102").Value = 1610.5264
$ws.Range("L102").Value = 1907.3334
$ws.Range("M102").Value = 11.47360000000003
$ws.Range("N102").Value = -5151.3334

# Row 131 (ARM)
$ws.Range("H131").Value = 78725
$ws.Range("J131").Value = 78725
$ws.Range("L131").Value = 78725
$ws.Range("N131").Value = -88805

# Row 132 (ARM)
$ws.Range("H132").Value = 1311.8372
$ws.Range("I132").Value = 1060.275
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 3180.825
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -650.8250000000003
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("CRP")
# Row 20 (CRP)
$ws.Range("H20").Value = 44467.6
$ws.Range("J20").Value = 44467.6
$ws.Range("L20").Value = 44467.6
$ws.Range("N20").Value = -44939.6

# Row 30 (CRP)
$ws.Range("H30").Value = 44467.6
$ws.Range("J30").Value = 44467.6
$ws.Range("L30").Value = 44467.6
$ws.Range("N30").Value = -44649.6

# Row 31 (CRP)
$ws.Range("H31").Value = 1905.2456
$ws.Range("I31").Value = 1425.4902
$ws.Range("J31").Value = 5983.1665
$ws.Range("K31").Value = 1425.4902
$ws.Range("L31").Value = 5983.1665
$ws.Range("M31").Value = -1130.4902
$ws.Range("N31").Value = -6573.1665

# Row 34 (CRP)
$ws.Range("H34").Value = 1905.2456
$ws.Range("I34").Value = 1425.4902
$ws.Range("J34").Value = 5983.1665
$ws.Range("K34").Value = 1425.4902
$ws.Range("L34").Value = 5983.1665
$ws.Range("M34").Value = -1223.4902
$ws.Range("N34").Value = -6387.1665

# Row 41 (CRP)
$ws.Range("H41").Value = 27786.334
$ws.Range("I41").Value = 7479.5
$ws.Range("J41").Value = 68400
$ws.Range("K41").Value = 7479.5
$ws.Range("L41").Value = 68400
$ws.Range("M41").Value = -7051.5
$ws.Range("N41").Value = -69256

# Row 51 (CRP)
$ws.Range("H51").Value = 13000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 13000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 13000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -14472

# Row 60 (CRP)
$ws.Range("H60").Value = 5000
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -4489
$ws.Range("N60").ClearContents()

# Row 61 (CRP)
$ws.Range("H61").Value = 13000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 13000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 13000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -13696

# Row 99 (CRP)
$ws.Range("H99").Value = 2727.5789
$ws.Range("I99").Value = 2289.1428
$ws.Range("J99").Value = 2983.3333
$ws.Range("K99").Value = 2289.1428
$ws.Range("L99").Value = 2983.3333
$ws.Range("M99").Value = -791.1428000000001
$ws.Range("N99").Value = -5979.3333

# Row 105 (CRP)
$ws.Range("H105").Value = 1672.8572
$ws.Range("I105").Value = 1402.5
$ws.Range("J105").Value = 2033.3334
$ws.Range("K105").Value = 1402.5
$ws.Range("L105").Value = 2033.3334
$ws.Range("M105").Value = 344.5
$ws.Range("N105").Value = -5527.3334

# Row 126 (CRP)
$ws.Range("H126").Value = 2727.5789
$ws.Range("I126").Value = 2289.1428
$ws.Range("J126").Value = 2983.3333
$ws.Range("K126").Value = 6867.428400000001
$ws.Range("L126").Value = 8949.999899999999
$ws.Range("M126").Value = -4397.428400000001
$ws.Range("N126").Value = -13889.9999

# Row 127 (CRP)
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

# Row 128 (CRP)
$ws.Range("H128").Value = 44467.6
$ws.Range("J128").Value = 44467.6
$ws.Range("L128").Value = 44467.6
$ws.Range("N128").Value = -54427.6

# Row 130 (CRP)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 109 (CUL)
$ws.Range("H109").Value = 2960.8696
$ws.Range("I109").Value = 415.3846
$ws.Range("J109").Value = 6270
$ws.Range("K109").Value = 1246.1538
$ws.Range("L109").Value = 18810
$ws.Range("M109").Value = -206.1538
$ws.Range("N109").Value = -20890

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (GSM)
$ws.Range("H132").Value = 11036.272
$ws.Range("I132").Value = 17233.5
$ws.Range("J132").Value = 3599.6
$ws.Range("K132").Value = 51700.5
$ws.Range("L132").Value = 10798.8
$ws.Range("M132").Value = -49170.5
$ws.Range("N132").Value = -15858.8

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 1000.5
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1002
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1002
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1378

# Row 100 (LTW)
$ws.Range("H100").Value = 2975
$ws.Range("I100").Value = 2800
$ws.Range("J100").Value = 3033.3333
$ws.Range("K100").Value = 2800
$ws.Range("L100").Value = 3033.3333
$ws.Range("M100").Value = -2259
$ws.Range("N100").Value = -4115.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Range("H62").Value = 3600.375
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3686.1428
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3686.1428
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4934.1428

# Row 65 (WVR)
$ws.Range("H65").Value = 3600.375
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3686.1428
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 18430.714
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -24670.714
